$d = $word.ActiveDocument

# -------------------------------------------------------------------
# 1) GPA: 3.11/4.0 -> GPA: 3.13/4.0, split across 3 runs with the
#    "_GoBack" bookmark relocated between the "3" and "/4.0" pieces.
# -------------------------------------------------------------------
$gpaRange = $d.Content
$found = $gpaRange.Find.Execute("GPA: 3.11/4.0")
if (-not $found) {
    throw "GPA text not found"
}
$gpaStart = $gpaRange.Start

# Force a run split between "GPA: 3.1" and the trailing "1/4.0" by
# dropping a throw-away bookmark at that boundary first.
$splitRange = $d.Range($gpaStart + 8, $gpaStart + 8)
$d.Bookmarks.Add("zzz_temp_split", $splitRange) | Out-Null

# Change the second "1" into a "3" (still inside its own run).
$digitRange = $d.Range($gpaStart + 8, $gpaStart + 9)
$digitRange.Text = "3"

# Relocate the (single, special) "_GoBack" bookmark to sit right after
# the new "3", before "/4.0".
$goBackRange = $d.Range($gpaStart + 9, $gpaStart + 9)
$d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null

# Drop the scaffolding bookmark now that the run boundary exists.
$d.Bookmarks("zzz_temp_split").Delete()

# -------------------------------------------------------------------
# 2) Merge the two runs making up the Heroku link back into one run.
# -------------------------------------------------------------------
$heroRange = $d.Content
$found2 = $heroRange.Find.Execute("(See https://csc309-team19.herokuapp.com)")
if (-not $found2) {
    throw "Heroku link text not found"
}
$heroScoped = $d.Range($heroRange.Start, $heroRange.End)
$heroScoped.Find.Execute(")", $true, $false, $false, $false, $false, $true, 1, $false, ")", 1) | Out-Null
